$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set F22 to a static value (removing its SUM formula)
$ws.Range("F22").Value = 1234

# Select cell I32 (matches updated <selection> in sheet view)
$ws.Range("I32").Select()
